$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 78. This shifts the existing
# rows 78-94 down to 79-95 (and the sheet dimension grows from R94 to R95).
$ws.Rows(78).Insert()

# The newly inserted row 78 is a new weekly price record for the same
# market/product/category as the row that used to be at 78 (now at 79);
# copy that row's formatting/static columns down into the blank row 78,
# then overwrite the fields that actually differ for the new record
# (date, volume, min/max/avg price, $/kg).
$ws.Range("A79:R79").Copy()
$ws.Range("A78:R78").PasteSpecial(-4104)
$excel.CutCopyMode = 0

$ws.Range("D78").Value2 = 44491
$ws.Range("J78").Value2 = 200
$ws.Range("K78").Value2 = 800
$ws.Range("L78").Value2 = 900
$ws.Range("M78").Value2 = 850
$ws.Range("P78").Value2 = 850
